$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 18474
$ws.Range("E3").Value = 1713
$ws.Range("E4").Value = 11700
$ws.Range("E5").Value = 17902
$ws.Range("E6").Value = 15421
$ws.Range("E7").Value = 1483
$ws.Range("E8").Value = 5235
$ws.Range("E9").Value = 13898
$ws.Range("E10").Value = 19552
$ws.Range("E11").Value = 5269
$ws.Range("E12").Value = 17354
$ws.Range("E13").Value = 5648
